$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New applicant row (row 3) appended below the existing header (row 1)
# and first applicant (row 2).
$ws.Range("A3").Value = "JOHAN"
$ws.Range("B3").Value = "ESTEBAN"
$ws.Range("C3").Value = "C"
$ws.Range("D3").Value = "C.C."
$ws.Range("F3").Value = "DIAGONAL AV-32"
$ws.Range("H3").Value = 18
$ws.Range("I3").Value = "LAURELES"
$ws.Range("J3").Value = 6
$ws.Range("L3").Value = "juanp@gmail.com"
$ws.Range("M3").Value = "CHIMBOALHOMBRO"
$ws.Range("N3").Value = 44
$ws.Range("O3").Value = "uploaded_files\132_CÉDULA.pdf"
$ws.Range("P3").Value = "uploaded_files\132_CIVICA.pdf"
$ws.Range("Q3").Value = "uploaded_files\132_SERVICIOPUBLICOS.pdf"
$ws.Range("R3").Value = "uploaded_files\132_ANEXO1.pdf"
$ws.Range("S3").Value = "uploaded_files\132_ANEXO2.xlsx"

# These look numeric/date-like but must stay as literal text (document
# number, birth date string, phone number), so force text format before
# assigning the value to avoid Excel auto-converting them.
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "132"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "09-01-2007"
$ws.Range("K3").NumberFormat = "@"
$ws.Range("K3").Value = "3002991878"
